$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.646.36'
$ws.Range('E2').Value = '  +5.33%  '

$ws.Range('D3').Value = '2.236.74'
$ws.Range('E3').Value = '  +3.49%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.82%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.87'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.66%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.403'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.26%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.81'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.43%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0878'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.21%  '

$ws.Range('E12').Value = '  +0.28%  '

$ws.Range('D13').Value = '2.569.08'
$ws.Range('E13').Value = '  +3.53%  '

$ws.Range('E14').Value = '  -1.86%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.800'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.98%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.57'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.34%  '

$ws.Range('D18').Value = '2.227.21'
$ws.Range('E18').Value = '  +3.19%  '

$ws.Range('D19').Value = '41.533.62'
$ws.Range('E19').Value = '  +4.63%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.47%  '

$ws.Range('D21').Value = '0.0₃0903'
$ws.Range('E21').Value = '  +6.78%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.25%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '247.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.51%  '

$ws.Range('E24').Value = '  +0.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.17%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.42%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.144'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.85%  '

$ws.Range('E29').Value = '  -2.10%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.09'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.24%  '

$ws.Range('E31').Value = '  +0.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.80'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.49%  '

$ws.Range('E33').Value = '  -0.47%  '

$ws.Range('E34').Value = '  +6.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.73%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0625'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.16%  '

$ws.Range('B37').Value = 'THORChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.71'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.43%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.79%  '

$ws.Range('E39').Value = '  -1.94%  '

$ws.Range('E40').Value = '  +0.52%  '

$ws.Range('B41').Value = 'FTXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.91'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.10%  '

$ws.Range('B42').Value = 'TerraClassic'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000237'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +24.10%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0237'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.33%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.82'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +14.41%  '

$ws.Range('E45').Value = '  -2.52%  '

$ws.Range('E46').Value = '  +3.71%  '

$ws.Range('D47').Value = '1.484.83'
$ws.Range('E47').Value = '  -2.78%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.18'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.29%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.66'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.64%  '

$ws.Range('E50').Value = '  -1.12%  '

$ws.Range('E51').Value = '  -2.13%  '
